$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Expected Behaviour" (validation) cell (H3) for VT200-0681: drop the
# validate5/Screenshot block that is no longer produced by the steps below.
$newValidate = "validate1`n{`nvalidate_PageTitle=Manual Compliance Ruby Specs`n};`nvalidate2`n{`nvalidate_PageTitle=Notification Ruby Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0681`n};"
$ws.Range("H3").Value = $newValidate

# Update the "Steps" cell (G3) for test case VT200-0681: drop the screenshot/validate5
# steps and replace the trailing wait/screenshot with a hide-popup text check.
$newSteps = "wait(5);`nvalidate1;`nlink_Click(notification_test_link);`nvalidate2;`nSelectTestToRun(VT200_0681_string);`nClickRunTest(runtest_top_xpath);`nwait(2);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nCheckUITextContains(This_is_a_pop_up_for_hide);`nwait(10);`nCheckUITextContains(hidepopup);"
$ws.Range("G3").Value = $newSteps

# The shorter text no longer needs as much wrapped height.
$ws.Rows.Item(3).RowHeight = 153.75

# Move the active selection to G2.
$ws.Range("G2").Select()
